# Section_4 WIP, finished upto lesson# 31 - Holt-Winters
# Log a new entry (row 14) in the time-tracking sheet and note the progress.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New session: date + start/end time (duration is computed by formula,
# matching the pattern used by the other rows in the log).
$ws.Range("B14").Value = 44826
$ws.Range("C14").Value = 0.9375
$ws.Range("D14").Value = 0.98958333333333337
$ws.Range("E14").Formula = "=D14-C14"

# Copy the "note" cell formatting (used for the progress remarks in column G)
# from the previous note at G13 onto G14, then write the new progress note.
$ws.Range("G13").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value = "Finished upto lesson# 31"

# Nudge the frozen-pane scroll position down to keep the newest rows in view.
$ws.Activate()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 7
}

$excel.CutCopyMode = 0
